$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 250.38889
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 729.8
$ws.Range("K2").Value = 66
$ws.Range("L2").Value = 729.8
$ws.Range("M2").Value = 47
$ws.Range("N2").Value = -955.8
$ws.Range("H29").Value = 123.833336
$ws.Range("I29").Value = 123.833336
$ws.Range("K29").Value = 371.500008
$ws.Range("M29").Value = -90.50000799999998
$ws.Range("H33").Value = 267
$ws.Range("I33").Value = 267
$ws.Range("K33").Value = 267
$ws.Range("M33").Value = -38
$ws.Range("H38").Value = 191.91667
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 4893.2383
$ws.Range("I40").Value = 3963.0908
$ws.Range("K40").Value = 3963.0908
$ws.Range("M40").Value = -3788.0908
$ws.Range("H41").Value = 492
$ws.Range("I41").Value = 492
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 492
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -52
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value = 1617.7333
$ws.Range("J58").Value = 2644.4443
$ws.Range("L58").Value = 7933.3329
$ws.Range("N58").Value = -8233.332900000001
$ws.Range("H86").Value = 1219
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1219
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1219
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -3465
$ws.Range("H89").Value = 1219
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1219
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 6095
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -17327
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("H111").Value = 1825.1515
$ws.Range("I111").Value = 1747.5
$ws.Range("K111").Value = 5242.5
$ws.Range("M111").Value = -2175.5
$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 1499
$ws.Range("K113").Value = 1499
$ws.Range("M113").Value = 1755
$ws.Range("H116").Value = 6023.077
$ws.Range("I116").Value = 5452.5
$ws.Range("K116").Value = 5452.5
$ws.Range("M116").Value = -2010.5
$ws.Range("H138").Value = 6645.6924
$ws.Range("J138").Value = 9999.714
$ws.Range("L138").Value = 29999.142
$ws.Range("N138").Value = -40279.142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10026398
$ws.Range("I32").Value = 23996.834
$ws.Range("J32").Value = 25030000
$ws.Range("K32").Value = 23996.834
$ws.Range("L32").Value = 25030000
$ws.Range("M32").Value = -23709.834
$ws.Range("N32").Value = -25030574
$ws.Range("H132").Value = 1722
$ws.Range("I132").Value = 1722
$ws.Range("K132").Value = 5166
$ws.Range("M132").Value = -2636
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H51").Value = 28187.5
$ws.Range("J51").Value = 42375
$ws.Range("L51").Value = 42375
$ws.Range("N51").Value = -43847
$ws.Range("H60").Value = 9377
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21022
$ws.Range("H61").Value = 28187.5
$ws.Range("J61").Value = 42375
$ws.Range("L61").Value = 42375
$ws.Range("N61").Value = -43071
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H141").Value = 52646.168
$ws.Range("J141").Value = 57175.4
$ws.Range("L141").Value = 57175.4
$ws.Range("N141").Value = -67535.39999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 94.92308
$ws.Range("I33").Value = 101.75
$ws.Range("J33").Value = 84
$ws.Range("K33").Value = 610.5
$ws.Range("L33").Value = 504
$ws.Range("M33").Value = -327.5
$ws.Range("N33").Value = -1070
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 6091.2
$ws.Range("I99").Value = 3485.3333
$ws.Range("K99").Value = 10455.9999
$ws.Range("M99").Value = -8209.999899999999
$ws.Range("H131").Value = 2094
$ws.Range("J131").Value = 1988
$ws.Range("L131").Value = 5964
$ws.Range("N131").Value = -16044
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 34376.7
$ws.Range("J57").Value = 46538.145
$ws.Range("L57").Value = 46538.145
$ws.Range("N57").Value = -48178.145
$ws.Range("H113").Value = 4177.826
$ws.Range("I113").Value = 2318.8572
$ws.Range("J113").Value = 7069.5557
$ws.Range("K113").Value = 2318.8572
$ws.Range("L113").Value = 7069.5557
$ws.Range("M113").Value = -148.8571999999999
$ws.Range("N113").Value = -11409.5557
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4596.5
$ws.Range("I7").Value = 2672.7144
$ws.Range("K7").Value = 2672.7144
$ws.Range("M7").Value = -2560.7144
$ws.Range("H68").Value = 3985.7144
$ws.Range("I68").Value = 1580
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 1580
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -831
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 3985.7144
$ws.Range("I71").Value = 1580
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 7900
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -4156
$ws.Range("N71").Value = -57488
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19002
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55008
$ws.Range("H82").Value = 2200
$ws.Range("I82").Value = 2200
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2200
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1839
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2200
$ws.Range("I85").Value = 2200
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2200
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -952
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 4596.5
$ws.Range("I126").Value = 2672.7144
$ws.Range("K126").Value = 8018.1432
$ws.Range("M126").Value = -5548.1432
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 768.4375
$ws.Range("I100").Value = 730.46155
$ws.Range("J100").Value = 933
$ws.Range("K100").Value = 1460.9231
$ws.Range("L100").Value = 1866
$ws.Range("M100").Value = -919.9231
$ws.Range("N100").Value = -2948
$ws.Range("H126").Value = 3851.36
$ws.Range("I126").Value = 3104.9412
$ws.Range("J126").Value = 5437.5
$ws.Range("K126").Value = 9314.8236
$ws.Range("L126").Value = 16312.5
$ws.Range("M126").Value = -6844.8236
$ws.Range("N126").Value = -21252.5
$ws.Range("H132").Value = 2870.5
$ws.Range("I132").Value = 2870.5
$ws.Range("K132").Value = 8611.5
$ws.Range("M132").Value = -6081.5
